# Auto-generated edit script: apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.435.65'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '1.792.14'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '226.25'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('E6').Value = '  +1.78%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '32.66'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.44%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.297'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.35%  '
$ws.Range('E10').Value = '  +0.62%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0950'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('D12').Value = '2.051.40'
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').Value = '1.829.13'
$ws.Range('E13').Value = '  +1.98%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '11.07'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('E15').Value = '  +1.96%  '
$ws.Range('D16').Value = '34.419.63'
$ws.Range('E16').Value = '  +0.80%  '
$ws.Range('E17').Value = '  +2.12%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '68.82'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '247.04'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('E20').Value = '  +2.82%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.21'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.58%  '
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.16'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.47%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.06'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.54%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '164.58'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.14%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.23'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.80%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.52'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.13%  '
$ws.Range('E28').Value = '  +2.62%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  +3.59%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('E32').Value = '  +0.47%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.89'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +7.14%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.82'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.09%  '
$ws.Range('D35').Value = '1.426.12'
$ws.Range('E35').Value = '  -1.06%  '
$ws.Range('E36').Value = '  +6.92%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.668'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +3.27%  '
$ws.Range('E38').Value = '  +2.00%  '
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '84.78'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +5.35%  '
$ws.Range('E41').Value = '  +1.17%  '
$ws.Range('E42').Value = '  +1.59%  '
$ws.Range('E43').Value = '  +1.74%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.59'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.75%  '
$ws.Range('E46').Value = '  +0.52%  '
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').Value = '1.948.48'
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '105.48'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.00'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0129'
$ws.Range('E51').Value = '  -4.78%  '
